$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header in H1, matching the formatting of the existing header cells (e.g. G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122) # xlPasteFormats
$ws.Application.CutCopyMode = $false
$ws.Range("H1").Value = "Save"

# New data value for row 2
$ws.Range("H2").Value = 0
